$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "P0102"
$ws.Range("B2").Value = "PR0001"
$ws.Range("C2").Value = "Resultados Carografía Social Taller San Marcos.pptx"
$ws.Range("D2").Value = "s3://foa-prod-comp-fenomenologico-bucket/foa_puj_curada/P0102/03_OTROS/Resultados Carograf├нa Social Taller San Marcos.pptx"
$ws.Range("E2").Value = 2024
$ws.Range("F2").Value = "Ministerio de Ambiente y Desarrollo Sostenible "
$ws.Range("K2").Value = "Es una presentación que resume y muestra cartográficamente los hallazgos de un taller de cartografía social realizado en San Marcos, en el 2024. En las diapositivas se presentan los puntos críticos que pueden ser afectados por inundaciones en siete municipios de La Mojana: Ayapel, Guaranda, Majagual, San Benito Abad, San Jacinto del Cauca, San Marcos y Sucre. "
$ws.Range("L2").Value = "Presentación"
$ws.Range("M2").Value = "pptx"
$ws.Range("N2").Value = "A partir del ejercicio de cartografía social, se identificaron causas de algunas afectaciones y las propuestas para abordarlas. Algunas de las afectaciones identificadas fueron la contaminación del agua, la salud, la pérdida de cultivos, de fauna y flora, el hacinamiento y las afectaciones emocionales y psicológicas. Como causas, se identificaron la minería, la construcción de terraplén sin tener en cuenta los cauces, la pérdida de capacidad hidráulica de los cauces. Algunas de las propuestas para afrontar estas situaciones fueron: capacitaciones a las comunidades, dragado de caños y ciénagas, control de la minería, compuertas en cierres y rompimientos, entre otros. Por último, se presentan los datos históricos desde 1950 de sequías e inundaciones en La Mojana."
$ws.Range("O2").Value = "Si"
$ws.Range("P2").Value = "s3://foa-prod-comp-fenomenologico-bucket/foa_puj_curada/P0102/03_OTROS\20240527_135418.jpg`ns3://foa-prod-comp-fenomenologico-bucket/foa_puj_curada/P0102/03_OTROS\20240527_135736.jpg`ns3://foa-prod-comp-fenomenologico-bucket/foa_puj_curada/P0102/03_OTROS\20240527_140057.jpg`ns3://foa-prod-comp-fenomenologico-bucket/foa_puj_curada/P0102/03_OTROS\Cartograf├нa social.xlsx`ns3://foa-prod-comp-fenomenologico-bucket/foa_puj_curada/P0102/03_OTROS\Matriz consolidada.xlsx"
$ws.Range("Q2").Value = "Variable: inundaciones"
$ws.Range("R2").Value = "Inundación, sequía, Mojana, afectaciones, minería, Ayapel, Guaranda, Majagual, San Benito Abad, San Jacinto del Cauca, San Marcos, Sucre."
